$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H80").Value = 517.7083
$ws_ALC.Range("I80").Value = 382.2
$ws_ALC.Range("J80").Value = 743.55554
$ws_ALC.Range("K80").Value = 1146.6
$ws_ALC.Range("L80").Value = 2230.66662
$ws_ALC.Range("M80").Value = -148.5999999999999
$ws_ALC.Range("N80").Value = -4226.66662

$ws_ALC.Range("H83").Value = 517.7083
$ws_ALC.Range("I83").Value = 382.2
$ws_ALC.Range("J83").Value = 743.55554
$ws_ALC.Range("K83").Value = 3439.8
$ws_ALC.Range("L83").Value = 6691.99986
$ws_ALC.Range("M83").Value = 1552.2
$ws_ALC.Range("N83").Value = -16675.99986

$ws_ALC.Range("H100").Value = 1547.6
$ws_ALC.Range("I100").Value = 1428.1
$ws_ALC.Range("K100").Value = 1428.1
$ws_ALC.Range("M100").Value = -887.0999999999999

$ws_ALC.Range("H125").Value = 980.7273
$ws_ALC.Range("I125").Value = 978.6667
$ws_ALC.Range("K125").Value = 8808.0003
$ws_ALC.Range("M125").Value = -6348.0003

$ws_ALC.Range("H132").Value = 2226.889
$ws_ALC.Range("I132").Value = 2277.4285
$ws_ALC.Range("K132").Value = 6832.2855
$ws_ALC.Range("M132").Value = -4302.2855

$ws_ALC.Range("H138").Value = 5081.477
$ws_ALC.Range("I138").Value = 2885.879
$ws_ALC.Range("K138").Value = 8657.636999999999
$ws_ALC.Range("M138").Value = -3517.636999999999

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H45").Value = 2721.182
$ws_ARM.Range("I45").Value = 2743.3
$ws_ARM.Range("K45").Value = 2743.3
$ws_ARM.Range("M45").Value = -2366.3

$ws_ARM.Range("H74").Value = 1123.8462
$ws_ARM.Range("J74").Value = 1499.5
$ws_ARM.Range("L74").Value = 1499.5
$ws_ARM.Range("N74").Value = -3247.5

$ws_ARM.Range("H77").Value = 1123.8462
$ws_ARM.Range("J77").Value = 1499.5
$ws_ARM.Range("L77").Value = 7497.5
$ws_ARM.Range("N77").Value = -16233.5

$ws_ARM.Range("H124").Value = 50000
$ws_ARM.Range("J124").Value = 50000
$ws_ARM.Range("L124").Value = 50000
$ws_ARM.Range("N124").Value = -59820

$ws_ARM.Range("H132").Value = 1754
$ws_ARM.Range("I132").Value = 1839.6666
$ws_ARM.Range("J132").Value = 1497
$ws_ARM.Range("K132").Value = 5518.9998
$ws_ARM.Range("L132").Value = 4491
$ws_ARM.Range("M132").Value = -2988.9998
$ws_ARM.Range("N132").Value = -9551

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 6749.5
$ws_BSM.Range("I20").Value = 11499
$ws_BSM.Range("K20").Value = 11499
$ws_BSM.Range("M20").Value = -11252

$ws_BSM.Range("H22").Value = 442.22223
$ws_BSM.Range("I22").Value = 376.42856
$ws_BSM.Range("J22").Value = 672.5
$ws_BSM.Range("K22").Value = 376.42856
$ws_BSM.Range("L22").Value = 672.5
$ws_BSM.Range("M22").Value = -203.42856
$ws_BSM.Range("N22").Value = -1018.5

$ws_BSM.Range("H44").Value = 25000
$ws_BSM.Range("J44").Value = 25000
$ws_BSM.Range("L44").Value = 25000
$ws_BSM.Range("N44").Value = -25994

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 2758.5557
$ws_CRP.Range("I31").Value = 1961.5238
$ws_CRP.Range("K31").Value = 1961.5238
$ws_CRP.Range("M31").Value = -1666.5238

$ws_CRP.Range("H34").Value = 2758.5557
$ws_CRP.Range("I34").Value = 1961.5238
$ws_CRP.Range("K34").Value = 1961.5238
$ws_CRP.Range("M34").Value = -1759.5238

$ws_CRP.Range("H58").Value = 2088.3572
$ws_CRP.Range("I58").Value = 1136.6364
$ws_CRP.Range("J58").Value = 5578
$ws_CRP.Range("K58").Value = 1136.6364
$ws_CRP.Range("L58").Value = 5578
$ws_CRP.Range("M58").Value = -933.6364000000001
$ws_CRP.Range("N58").Value = -5984

$ws_CRP.Range("H99").Value = 14255.19
$ws_CRP.Range("I99").Value = 12334.143
$ws_CRP.Range("J99").Value = 15215.714
$ws_CRP.Range("K99").Value = 12334.143
$ws_CRP.Range("L99").Value = 15215.714
$ws_CRP.Range("M99").Value = -10836.143
$ws_CRP.Range("N99").Value = -18211.714

$ws_CRP.Range("H124").Value = 99999
$ws_CRP.Range("J124").Value = 99999
$ws_CRP.Range("L124").Value = 99999
$ws_CRP.Range("N124").Value = -104909

$ws_CRP.Range("H126").Value = 14255.19
$ws_CRP.Range("I126").Value = 12334.143
$ws_CRP.Range("J126").Value = 15215.714
$ws_CRP.Range("K126").Value = 37002.429
$ws_CRP.Range("L126").Value = 45647.142
$ws_CRP.Range("M126").Value = -34532.429
$ws_CRP.Range("N126").Value = -50587.142

$ws_CRP.Range("H134").Value = 3722.5
$ws_CRP.Range("I134").Value = 3596.2856
$ws_CRP.Range("J134").Value = 3899.2
$ws_CRP.Range("K134").Value = 10788.8568
$ws_CRP.Range("L134").Value = 11697.6
$ws_CRP.Range("M134").Value = -8253.856800000001
$ws_CRP.Range("N134").Value = -16767.6

$ws_CRP.Range("H136").Value = 2088.3572
$ws_CRP.Range("I136").Value = 1136.6364
$ws_CRP.Range("J136").Value = 5578
$ws_CRP.Range("K136").Value = 3409.9092
$ws_CRP.Range("L136").Value = 16734
$ws_CRP.Range("M136").Value = -859.9092000000001
$ws_CRP.Range("N136").Value = -21834

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H26").Value = 0
$ws_CUL.Range("J26").Value = 0
$ws_CUL.Range("L26").Value = 0
$ws_CUL.Range("N26").Value = $null

$ws_CUL.Range("H137").Value = 3855.1875
$ws_CUL.Range("I137").Value = 3366.3333
$ws_CUL.Range("J137").Value = 4148.5
$ws_CUL.Range("K137").Value = 10098.9999
$ws_CUL.Range("L137").Value = 12445.5
$ws_CUL.Range("M137").Value = -4998.999899999999
$ws_CUL.Range("N137").Value = -22645.5

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H26").Value = 15000
$ws_GSM.Range("J26").Value = 15000
$ws_GSM.Range("L26").Value = 15000
$ws_GSM.Range("N26").Value = -15560

$ws_GSM.Range("H50").Value = 15000
$ws_GSM.Range("J50").Value = 15000
$ws_GSM.Range("L50").Value = 15000
$ws_GSM.Range("N50").Value = -15996

$ws_GSM.Range("H70").Value = 6199.6
$ws_GSM.Range("I70").Value = 5334.5
$ws_GSM.Range("K70").Value = 5334.5
$ws_GSM.Range("M70").Value = -5064.5

$ws_GSM.Range("H73").Value = 6199.6
$ws_GSM.Range("I73").Value = 5334.5
$ws_GSM.Range("K73").Value = 5334.5
$ws_GSM.Range("M73").Value = -4398.5

$ws_GSM.Range("H97").Value = 388
$ws_GSM.Range("I97").Value = 401.2857
$ws_GSM.Range("J97").Value = 295
$ws_GSM.Range("K97").Value = 401.2857
$ws_GSM.Range("L97").Value = 295
$ws_GSM.Range("M97").Value = 94.71429999999998
$ws_GSM.Range("N97").Value = -1287

$ws_GSM.Range("H113").Value = 3129.182
$ws_GSM.Range("I113").Value = 2801
$ws_GSM.Range("J113").Value = 4606
$ws_GSM.Range("K113").Value = 2801
$ws_GSM.Range("L113").Value = 4606
$ws_GSM.Range("M113").Value = -631
$ws_GSM.Range("N113").Value = -8946

$ws_GSM.Range("H122").Value = 61568.35
$ws_GSM.Range("I122").Value = 2057.9167
$ws_GSM.Range("K122").Value = 6173.750100000001
$ws_GSM.Range("M122").Value = -3723.750100000001

$ws_GSM.Range("H132").Value = 2731
$ws_GSM.Range("I132").Value = 2406.8
$ws_GSM.Range("J132").Value = 3271.3333
$ws_GSM.Range("K132").Value = 7220.400000000001
$ws_GSM.Range("L132").Value = 9813.999899999999
$ws_GSM.Range("M132").Value = -4690.400000000001
$ws_GSM.Range("N132").Value = -14873.9999

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H132").Value = 2544.8845
$ws_LTW.Range("I132").Value = 2198.8572
$ws_LTW.Range("K132").Value = 6596.571599999999
$ws_LTW.Range("M132").Value = -4066.571599999999

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H51").Value = 1967
$ws_WVR.Range("I51").Value = 1967
$ws_WVR.Range("K51").Value = 1967
$ws_WVR.Range("M51").Value = -1457

$ws_WVR.Range("H52").Value = 15000
$ws_WVR.Range("I52").Value = 0
$ws_WVR.Range("J52").Value = 15000
$ws_WVR.Range("K52").Value = 0
$ws_WVR.Range("L52").Value = 15000
$ws_WVR.Range("M52").Value = $null
$ws_WVR.Range("N52").Value = -15452

$ws_WVR.Range("H54").Value = 35000
$ws_WVR.Range("I54").Value = 0
$ws_WVR.Range("J54").Value = 35000
$ws_WVR.Range("K54").Value = 0
$ws_WVR.Range("L54").Value = 35000
$ws_WVR.Range("M54").Value = $null
$ws_WVR.Range("N54").Value = -36040

$ws_WVR.Range("H62").Value = 7857.2856
$ws_WVR.Range("J62").Value = 7857.2856
$ws_WVR.Range("L62").Value = 7857.2856
$ws_WVR.Range("N62").Value = -9105.285599999999

$ws_WVR.Range("H65").Value = 7857.2856
$ws_WVR.Range("J65").Value = 7857.2856
$ws_WVR.Range("L65").Value = 39286.428
$ws_WVR.Range("N65").Value = -45526.428

$ws_WVR.Range("H107").Value = 1898.75
$ws_WVR.Range("I107").Value = 1898.3334
$ws_WVR.Range("K107").Value = 5695.0002
$ws_WVR.Range("M107").Value = -3775.0002

$ws_WVR.Range("H132").Value = 3208.3704
$ws_WVR.Range("I132").Value = 2890.6316
$ws_WVR.Range("J132").Value = 3963
$ws_WVR.Range("K132").Value = 8671.8948
$ws_WVR.Range("L132").Value = 11889
$ws_WVR.Range("M132").Value = -6141.8948
$ws_WVR.Range("N132").Value = -16949
